# Applies the diff: adds two new worksheets (o_20, o_20_jumbled) after o_10,
# and adds a new "evaluator_partial_correctness" column (E) to all three sheets.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Create the two new worksheets, in order, right after o_10 ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "o_20"
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "o_20_jumbled"

# --- Re-select o_10 as the active / tabSelected sheet ---
$ws1.Activate()

# --- Sheet 1 (o_10): add new header cell (copy D1 format first so it gets the bold/border/center style) ---
$ws1.Range("D1").Copy($ws1.Range("E1"))
$ws1.Range("E1").Value = 'evaluator_partial_correctness'

# --- Sheet 1 (o_10): refresh row 2 values (prompt + llm_response changed, new evaluator_partial_correctness cell) ---
$ws1.Range("A2").Value = ' Given is the adjacency matrix for a unweighted undirected graph containing 10 nodes labelled A to J. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   

Consider some examples

Example 1: is the following a valid eulerian graph, if traversal is started from 0?
   A B C D E F G H I J K L M
 A 0 0 1 0 0 1 0 1 0 0 0 0 0
 B 0 0 0 0 0 0 0 0 1 0 0 1 0
 C 1 0 0 1 1 0 1 0 0 0 0 1 0
 D 0 0 1 0 0 1 0 0 0 0 0 0 0
 E 0 0 1 0 0 1 1 0 0 0 0 1 0
 F 1 0 0 1 1 0 0 0 1 0 0 1 0
 G 0 0 1 0 1 0 0 0 0 0 0 1 0
 H 1 0 0 0 0 0 0 0 1 0 1 0 0
 I 0 1 0 0 0 1 0 1 0 0 0 0 1
 J 0 0 0 0 0 0 0 0 0 0 1 0 0
 K 0 0 0 0 0 0 0 1 0 1 0 0 1
 L 0 1 1 0 1 1 1 0 0 0 0 0 1
 M 0 0 0 0 0 0 0 0 1 0 1 1 0

Solution: This is not a valid eulerian graph
        
 Given these examples, answer the following quesiton.

is the following a valid eulerian graph, if traversal is started from 1?

   A B C D E F G H I J
 A 0 1 1 0 0 0 0 0 0 1
 B 1 0 1 0 0 0 0 0 0 0
 C 1 1 0 1 0 1 0 0 1 0
 D 0 0 1 0 1 0 0 0 0 0
 E 0 0 0 1 0 1 1 0 0 1
 F 0 0 1 0 1 0 1 0 0 1
 G 0 0 0 0 1 1 0 1 0 1
 H 0 0 0 0 0 0 1 0 1 0
 I 0 0 1 0 0 0 0 1 0 0
 J 1 0 0 0 1 1 1 0 0 0
    '
$ws1.Range("B2").Value = 'This is a valid eulerian graph'
$ws1.Range("C2").Value = 'To determine whether the given graph is a valid eulerian graph, we can use the following criteria:
1. All nodes must have an even degree (i.e., the sum of connections of each node must be even). 
If all the nodes in the graph satisfy this criterion, then the graph is eulerian. If there are exactly two nodes with odd degrees, then the graph has an eulerian path, but not an eulerian circuit. If any other node has an odd degree, then the graph is not eulerian.
Let''s analyze the given graph:
A has a degree of 3, which is odd.
B has a degree of 2, which is even.
C has a degree of 5, which is odd.
D has a degree of 2, which is even.
E has a degree of 4, which is even.
F has a degree of 4, which is even.
G has a degree of 4, which is even.
H has a degree of 2, which is even.
I has a degree of 3, which is odd.
J has a degree of 4, which is even.
In this graph, A, C, and I have odd degrees, which means it is not a valid eulerian graph.
Therefore, the answer to the question is: No, this graph is not a valid eulerian graph if traversal is started from 1.'
$ws1.Range("D2").Value = 'Wrong'
$ws1.Range("E2").Value = 'N/A'
$ws1.Rows.Item(2).AutoFit()

# --- Copy the (now 5-column) header row formatting/values to the new sheets ---
$ws1.Range("A1:E1").Copy($ws2.Range("A1:E1"))
$ws1.Range("A1:E1").Copy($ws3.Range("A1:E1"))

# --- Sheet 2 (o_20) data row ---
$ws2.Range("A2").Value = ' Given is the adjacency matrix for a unweighted undirected graph containing 20 nodes labelled A to T. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   
Consider some examples
Example 1: is the following a valid eulerian graph, if traversal is started from 0?
   A B C D E F G H I J K L M N O P Q R S
 A 0 0 0 0 0 0 1 1 0 0 0 1 0 1 0 0 0 1 0
 B 0 0 0 0 0 0 0 0 0 0 0 0 0 1 1 0 0 0 0
 C 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0
 D 0 0 0 0 0 1 0 1 0 0 1 0 1 1 0 0 0 1 1
 E 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0
 F 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 1 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0
 H 1 0 0 1 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0
 I 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 1 1 0
 J 0 0 0 0 0 0 0 0 0 0 0 1 1 0 0 0 0 0 1
 K 0 0 1 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1
 L 1 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 1 0 0
 M 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0
 N 1 1 0 1 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1
 O 0 1 0 0 0 0 1 1 0 0 0 0 0 1 0 1 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0
 Q 0 0 0 0 0 0 0 0 1 0 0 1 0 0 0 0 0 0 0
 R 1 0 0 1 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 S 0 0 0 1 0 0 0 0 0 1 1 0 0 1 0 0 0 0 0
Solution: This is not a valid eulerian graph
 Given these examples, answer the following quesiton.
is the following a valid eulerian graph, if traversal is started from A?
   A B C D E F G H I J K L M N O P Q R S T
 A 0 0 0 1 0 0 0 0 0 0 0 0 0 0 1 0 0 1 1 0
 B 0 0 0 0 1 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0
 C 0 0 0 1 0 0 0 0 0 0 0 0 1 0 0 1 0 0 0 0
 D 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 1
 E 0 1 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 0 0 1 0 0 0 0 0 0 0 1 0 0 1 0 1
 G 0 1 0 1 0 1 0 0 0 0 0 0 0 0 1 0 0 0 0 0
 H 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 1 0 0 0 0
 I 0 1 0 0 1 0 0 0 0 1 0 0 0 0 1 0 0 1 0 1
 J 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 1 1 0 0
 K 0 0 0 0 0 0 0 1 0 0 0 1 0 0 1 0 0 0 0 1
 L 0 0 0 0 0 0 0 0 0 0 1 0 0 1 0 0 0 0 0 0
 M 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0
 O 1 0 0 0 0 1 1 0 1 0 1 0 0 0 0 0 1 0 0 0
 P 0 0 1 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 Q 0 0 0 1 0 0 0 0 0 1 0 0 0 0 1 0 0 0 0 0
 R 1 0 0 0 0 1 0 0 1 1 0 0 0 1 0 0 0 0 0 0
 S 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 T 0 0 0 1 0 1 0 0 1 0 1 0 0 0 0 0 0 0 0 0
    '
$ws2.Range("B2").Value = 'This is not a valid eulerian graph'
$ws2.Range("C2").Value = 'To determine if the graph is a valid eulerian graph, we need to check if each vertex has an even degree.
Counting the degree of each vertex:
- Vertex A has a degree of 3
- Vertex B has a degree of 4
- Vertex C has a degree of 3
- Vertex D has a degree of 5
- Vertex E has a degree of 2
- Vertex F has a degree of 4
- Vertex G has a degree of 4
- Vertex H has a degree of 2
- Vertex I has a degree of 7
- Vertex J has a degree of 3
- Vertex K has a degree of 4
- Vertex L has a degree of 2
- Vertex M has a degree of 1
- Vertex N has a degree of 2
- Vertex O has a degree of 6
- Vertex P has a degree of 2
- Vertex Q has a degree of 3
- Vertex R has a degree of 4
- Vertex S has a degree of 1
- Vertex T has a degree of 4
From the given degrees, we can see that Vertex M and Vertex S have odd degrees. Therefore, the graph is not a valid eulerian graph.'
$ws2.Range("D2").Value = 'Correct'
$ws2.Range("E2").Value = 'N/A'
$ws2.Rows.Item(2).AutoFit()

# --- Sheet 3 (o_20_jumbled) data row ---
$ws3.Range("A2").Value = ' Given is the adjacency matrix for a unweighted undirected graph containing 20 nodes labelled A to T. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   
Consider some examples
Example 1: is the following a valid eulerian graph, if traversal is started from 0?
   A B C D E F G H I J K L M N O P Q R S
 A 0 0 0 0 0 0 1 1 0 0 0 1 0 1 0 0 0 1 0
 B 0 0 0 0 0 0 0 0 0 0 0 0 0 1 1 0 0 0 0
 C 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0
 D 0 0 0 0 0 1 0 1 0 0 1 0 1 1 0 0 0 1 1
 E 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0
 F 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 1 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0
 H 1 0 0 1 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0
 I 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 1 1 0
 J 0 0 0 0 0 0 0 0 0 0 0 1 1 0 0 0 0 0 1
 K 0 0 1 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1
 L 1 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 1 0 0
 M 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0
 N 1 1 0 1 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1
 O 0 1 0 0 0 0 1 1 0 0 0 0 0 1 0 1 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0
 Q 0 0 0 0 0 0 0 0 1 0 0 1 0 0 0 0 0 0 0
 R 1 0 0 1 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 S 0 0 0 1 0 0 0 0 0 1 1 0 0 1 0 0 0 0 0
Solution: This is not a valid eulerian graph
 Given these examples, answer the following quesiton.
is the following a valid eulerian graph, if traversal is started from A?
   A B C D E F G H I J K L M N O P Q R S T
 A 0 1 1 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 B 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 1 0
 C 1 0 0 0 1 1 0 0 0 0 0 1 1 0 1 0 0 0 0 0
 D 0 0 0 0 0 0 0 1 0 0 0 1 1 0 0 0 0 0 1 0
 E 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 1 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0
 H 0 0 0 1 0 1 0 0 0 0 0 0 1 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0
 J 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 1 0 0
 K 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0
 L 0 0 1 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 M 0 0 1 1 0 0 0 1 0 0 0 0 0 0 0 1 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 1
 O 0 0 1 0 0 0 1 0 0 1 0 0 0 1 0 0 0 0 0 1
 P 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0
 Q 0 0 0 0 0 0 1 0 1 0 0 0 0 1 0 1 0 0 0 0
 R 0 1 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 S 0 1 0 1 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 1 1 0 0 0 0 0
    '
$ws3.Range("B2").Value = 'This is not a valid eulerian graph'
$ws3.Range("C2").Value = 'This is not a valid Eulerian graph if the traversal is started from A. In an Eulerian graph, every vertex must have an even degree. In this graph, vertices D, E, F, H, J, L, M, N, O, Q, S have an odd degree, which violates the necessary condition for an Eulerian graph.'
$ws3.Range("D2").Value = 'Correct'
$ws3.Range("E2").Value = 'N/A'
$ws3.Rows.Item(2).AutoFit()

